# M4_Seurat_CellType.pptx edits
# 1) Footer "Last updated" date field on the slide master / all 5 layouts /
#    notes master: 4/3/2024 -> 4/4/24
# 2) Slide 1 (title slide): title "Module 5" -> "Module 4", and the
#    subtitle text/format/position reworked.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholders - update the cached datetimeFigureOut text everywhere
#    it lives: the slide master, every slide layout, and the notes master.
# ---------------------------------------------------------------------------
function Set-DateShapeText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "4/4/24"
        }
    }
}

Set-DateShapeText $p.SlideMaster.Shapes

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Set-DateShapeText $layout.Shapes
}

Set-DateShapeText $p.NotesMaster.Shapes

# ---------------------------------------------------------------------------
# 2) Slide 1 - title + subtitle
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

$title = $s1.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Module 4: Cell Type Identification"

$sub = $s1.Shapes.Item(2)

# Move the subtitle placeholder down (explicit xfrm override).
$sub.Top = 321.25

# Re-point the text runs: "NCSU " -> "NC State ", keep "scRNA", and
# "-Seq Workshop, 2024" -> " Workshop, 2024".
$tr = $sub.TextFrame.TextRange
$run1 = $tr.Characters(1, 5)
$run1.Text = "NC State "

$tr = $sub.TextFrame.TextRange
$run3 = $tr.Characters(15, 19)
$run3.Text = " Workshop, 2024"

# Bump every run to 40pt.
$tr = $sub.TextFrame.TextRange
$tr.Font.Size = 40

# Force the body to not auto-shrink (matches <a:noAutofit/>).
$sub.TextFrame.AutoSize = 0
